$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the paragraph run that currently reads:
#   "link search and related to watched vids" + spaces + "////////////done"
# and rewrite it so "watched vids" becomes "history", the run is split into
# three runs (prefix / "history" / suffix) and the stretch of spaces before
# the slashes is shortened, matching the target diff.
# ---------------------------------------------------------------------------

$oldText = "link search and related to watched vids" + `
    "                                                              " + `
    "////////////done"

$newPart1 = "link search and related to "
$newPart2 = "history"
$newPart3 = "                                                      ////////////done"
$newFull  = $newPart1 + $newPart2 + $newPart3

# Find the exact old run text so we know precisely where it lives.
$rngFind = $d.Content
$found = $rngFind.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate target text to edit"
}

$startPos = $rngFind.Start

# Replace the whole stretch of text in one go (keeps a single run for now).
$rngFind.Text = $newFull

# Compute the boundaries of the three desired runs in absolute document
# character positions.
$r1Start = $startPos
$r1End   = $r1Start + $newPart1.Length
$r2Start = $r1End
$r2End   = $r2Start + $newPart2.Length
$r3Start = $r2End
$r3End   = $r3Start + $newPart3.Length

# Force Word to split the single run into three runs at the desired
# boundaries by nudging (and reverting) a character formatting attribute on
# each of the trailing pieces. The net formatting is unchanged, but the run
# boundary remains.
$rng2 = $d.Range($r2Start, $r2End)
$rng2.Font.Bold = 1
$rng2.Font.Bold = 0

$rng3 = $d.Range($r3Start, $r3End)
$rng3.Font.Bold = 1
$rng3.Font.Bold = 0
